# Apply the recorded edits to Statistics.xlsx (Sheet1):
#  - Swap the "RMS"/"Std" header labels in E2/F2
#  - Swap the RMS/Std measurement values in E/F for the affected data rows
#  - Update the saved view (selection, column width, default row height)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row: E2 <-> F2 (RMS/Std swap) ---
$e2 = $ws.Range("E2").Value()
$f2 = $ws.Range("F2").Value()
$ws.Range("E2").Value = $f2
$ws.Range("F2").Value = $e2

# --- Data rows whose E/F (RMS/Std) values need to be swapped ---
$rowsToSwap = @(33,34,35,36,37,38,39,40,41,43,47,48,49,50,51,52)
foreach ($r in $rowsToSwap) {
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $fCell = $ws.Cells.Item($r, 6)   # column F
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()
    $eCell.Value = $fVal
    $fCell.Value = $eVal
}

# --- View/formatting touch-ups ---
# (14.15 rounds to a stored column width of 15, matching the target width)
$ws.Columns("B").ColumnWidth = 14.15
$ws.Range("H5").Select()
